# Updated cryptos list on Mon Sep 23 21:54:25 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to stay a text string (the Price column stores
    # numbers formatted like "609.90" / "0.0240" as literal text, and
    # plain .Value assignment would let Excel auto-coerce these
    # number-looking strings into floating point values, losing the
    # exact original formatting).
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "63.319.14"
$ws.Range("E2").Value = "  +0.33%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.655.98"
$ws.Range("E3").Value = "  +3.50%  "

# Row 4 - TetherUSD (price unchanged)
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
Set-TextValue "D5" "609.90"
$ws.Range("E5").Value = "  +4.16%  "

# Row 6 - Solana
Set-TextValue "D6" "143.83"
$ws.Range("E6").Value = "  +0.00%  "

# Row 7 - USDC
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  -0.06%  "

# Row 8 - XRP (price unchanged)
$ws.Range("E8").Value = "  -0.41%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "2.653.15"
$ws.Range("E9").Value = "  +3.47%  "

# Row 10 - Dogecoin (price unchanged)
$ws.Range("E10").Value = "  +1.36%  "

# Row 11 - Toncoin (price unchanged)
$ws.Range("E11").Value = "  +0.01%  "

# Row 12 - TRON (price unchanged)
$ws.Range("E12").Value = "  +0.38%  "

# Row 13 - Cardano (price unchanged)
$ws.Range("E13").Value = "  +3.27%  "

# Row 14 - Avalanche
Set-TextValue "D14" "27.31"
$ws.Range("E14").Value = "  +1.04%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "3.128.60"
$ws.Range("E15").Value = "  +3.39%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "63.138.84"
$ws.Range("E16").Value = "  +0.19%  "

# Row 17 - ShibaInu (price unchanged)
$ws.Range("E17").Value = "  -0.39%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.689.67"
$ws.Range("E18").Value = "  +5.00%  "

# Row 19 - Chainlink
Set-TextValue "D19" "11.43"
$ws.Range("E19").Value = "  +3.82%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "341.46"
$ws.Range("E20").Value = "  +0.49%  "

# Row 21 - Polkadot (price unchanged)
$ws.Range("E21").Value = "  +2.26%  "

# Row 22 - Uniswap (price unchanged)
$ws.Range("E22").Value = "  +3.81%  "

# Row 24 - Litecoin
Set-TextValue "D24" "66.91"
$ws.Range("E24").Value = "  -1.08%  "

# Row 25 - Fetch.AI (price unchanged)
$ws.Range("E25").Value = "  +2.95%  "

# Row 26 - SuiNetwork (price unchanged)
$ws.Range("E26").Value = "  +0.26%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "8.67"
$ws.Range("E27").Value = "  +6.18%  "

# Row 28 - Kaspa (price unchanged)
$ws.Range("E28").Value = "  +0.44%  "

# Row 29 - Bittensor
Set-TextValue "D29" "545.94"
$ws.Range("E29").Value = "  +16.78%  "

# Row 30 - Binance-PegBSC-USD (price unchanged)
$ws.Range("E30").Value = "  -0.03%  "

# Row 31 - Aptos
Set-TextValue "D31" "7.81"
$ws.Range("E31").Value = "  -1.38%  "

# Row 33 - ImmutableX
Set-TextValue "D33" "1.78"
$ws.Range("E33").Value = "  +7.34%  "

# Row 35 - Monero
Set-TextValue "D35" "172.08"
$ws.Range("E35").Value = "  -2.15%  "

# Row 36 - NEARProtocol
Set-TextValue "D36" "5.12"
$ws.Range("E36").Value = "  +13.08%  "

# Row 37 - PolygonEcosystemToken (price unchanged)
$ws.Range("E37").Value = "  +2.64%  "

# Row 38 - FirstDigitalUSD (price unchanged)
$ws.Range("E38").Value = "  -0.12%  "

# Row 39 - EthereumClassic
Set-TextValue "D39" "19.11"
$ws.Range("E39").Value = "  +1.73%  "

# Row 40 - Stacks
Set-TextValue "D40" "1.86"
$ws.Range("E40").Value = "  +10.04%  "

# Row 41 - Aave
Set-TextValue "D41" "174.31"
$ws.Range("E41").Value = "  +11.05%  "

# Row 42 - USDe (price unchanged)
$ws.Range("E42").Value = "  -0.01%  "

# Row 43 - Filecoin (price unchanged)
$ws.Range("E43").Value = "  +2.10%  "

# Row 44 - InjectiveProtocol
Set-TextValue "D44" "22.21"
$ws.Range("E44").Value = "  +5.02%  "

# Row 45 - Hedera
Set-TextValue "D45" "0.0573"
$ws.Range("E45").Value = "  +6.80%  "

# Row 46 - Mantle
Set-TextValue "D46" "0.631"
$ws.Range("E46").Value = "  +0.16%  "

# Row 47 - Stellar
Set-TextValue "D47" "0.0962"
$ws.Range("E47").Value = "  +0.28%  "

# Row 48 - VeChain
Set-TextValue "D48" "0.0240"
$ws.Range("E48").Value = "  +1.57%  "

# Row 49 - EnergySwap
Set-TextValue "D49" "18.72"
$ws.Range("E49").Value = "  +4.10%  "

# Row 50 - dogwifhat (price unchanged)
$ws.Range("E50").Value = "  +5.21%  "

# Row 51 - WhiteBITCoin
Set-TextValue "D51" "11.25"
$ws.Range("E51").Value = "  -1.19%  "
